# This workbook lists Metric Impact-Resistant Polycarbonate Pan Head Phillips
# Screws grouped by thread size. Each group currently starts with a "banner"
# row that holds only the thread size (e.g. "M3 x 0.5 mm") in column A and
# "Polycarbonate Plastic" in column M (material_surface), with every other
# cell in that row blank.
#
# The edit:
#   1. Pushes each banner row's thread-size text down into column L
#      (thread_size) of every data row that belongs to that group, then
#      removes the now-redundant banner rows.
#   2. Replaces the old text header row (row 1: "Lg., mm", "Threading", ...)
#      with a new numeric index row (0, 1, 2, ... 12) - this becomes the new
#      row 1 - while the old header text becomes row 2 (with the former
#      "thread_size"/"material_surface" header labels in L2/M2 cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: propagate each group's thread-size banner into column L (12) of
# its member data rows, before the banner rows are removed. ---
$groups = @(
  @{Header = 2;  Rows = @(3, 4, 5, 6, 7)},
  @{Header = 8;  Rows = @(9, 10, 11, 12, 13)},
  @{Header = 14; Rows = @(15, 16, 17, 18, 19)},
  @{Header = 20; Rows = @(21, 22, 23, 24, 25)}
)

foreach ($g in $groups) {
  $label = $ws.Cells.Item($g.Header, 1).Value2
  foreach ($r in $g.Rows) {
    $ws.Cells.Item($r, 12).Value2 = $label
  }
}

# --- Step 2: delete the 4 banner rows, bottom-to-top so earlier row numbers
# stay valid as each delete shifts the rows below it up. ---
$ws.Rows.Item(20).Delete()
$ws.Rows.Item(14).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(2).Delete()

# --- Step 3: the header row's old "thread_size" / "material_surface" labels
# (L1/M1) no longer apply once it becomes row 2, so clear them. ---
$ws.Cells.Item(1, 12).Value2 = ""
$ws.Cells.Item(1, 13).Value2 = ""

# --- Step 4: insert a brand-new row 1 for the numeric index row. This shifts
# the current row 1 (the text header, along with its bold/border/center
# style) down to row 2. ---
$ws.Rows.Item(1).Insert()

# --- Step 5: restore correct styling.
# Copy the header formatting (now sitting on row 2 after the insert) back up
# onto the new row 1 so it reuses the existing bold/border/centered style,
# then reset row 2 back to the plain default "Normal" style. ---
$ws.Range("A2:M2").Copy()
$ws.Range("A1:M1").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("A2:M2").Style = "Normal"

# --- Step 6: populate the new row 1 with the sequential numeric index
# 0 .. 12 across columns A .. M. ---
for ($c = 1; $c -le 13; $c++) {
  $ws.Cells.Item(1, $c).Value2 = $c - 1
}
